# Test_cases.xlsx: mark Integration / Test gemaakt / Passed columns (F, I, J)
# with "X" for the test cases that have now passed, matching column K
# ("Werkt in app") which already had the marks. Also fill in rows 34, 35
# and 37 (TC17, TC18, TC19), which previously had no results at all.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows that already had an "X" in column K and now also get "X" in F, I, J
$rowsWithExistingResult = @(4, 6, 8, 11, 13, 16, 17, 19, 20, 21, 23, 24, 26, 27, 29, 31)

foreach ($r in $rowsWithExistingResult) {
    $ws.Range("F" + $r).Value = "X"
    $ws.Range("I" + $r).Value = "X"
    $ws.Range("J" + $r).Value = "X"
}

# Rows 34, 35 and 37 previously had no F/I/J/K values at all; fill them
# in now that those cases have been implemented and tested too. (The
# source data really does read "TC06"/"TC07"/"TC14" here rather than
# "X" -- matching the author's own workbook exactly.)
$ws.Range("F34").Value = "TC06"
$ws.Range("I34").Value = "TC06"
$ws.Range("J34").Value = "TC06"
$ws.Range("K34").Value = "X"

$ws.Range("F35").Value = "TC07"
$ws.Range("I35").Value = "TC07"
$ws.Range("J35").Value = "TC07"
$ws.Range("K35").Value = "X"

$ws.Range("F37").Value = "TC14"
$ws.Range("I37").Value = "TC14"
$ws.Range("J37").Value = "TC14"
$ws.Range("K37").Value = "X"

# Restore the view state (scroll position / selection) to where the
# author left it after finishing the test pass.
$ws.Activate()
$ws.Range("E31").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
